$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data per the Nov 14 2023 refresh.
# Column D (Price) values that look like plain numbers must be forced to
# remain text cells (matching the original inlineStr storage), otherwise
# Excel auto-converts them to numeric cells and mangles values such as
# "93.10" -> 93.1 or "4.80" -> 4.8.

$ws.Range("D2").Value = '36.454.93'
$ws.Range("E2").Value = '  -2.13%  '

$ws.Range("D3").Value = '2.053.49'
$ws.Range("E3").Value = '  -0.40%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.53'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -3.13%  '

$ws.Range("E6").Value = '  -0.61%  '

$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '54.22'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -8.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '57.99'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -4.52%  '

$ws.Range("E10").Value = '  -8.01%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0747'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -5.65%  '

$ws.Range("E12").Value = '  -2.59%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.895'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -2.88%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '14.58'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -8.71%  '

$ws.Range("D15").Value = '2.353.30'
$ws.Range("E15").Value = '  -0.46%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.33'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -8.96%  '

$ws.Range("D17").Value = '2.055.58'
$ws.Range("E17").Value = '  -0.95%  '

$ws.Range("D18").Value = '36.374.86'
$ws.Range("E18").Value = '  -2.34%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.58'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -12.07%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.84'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -4.88%  '

$ws.Range("D21").Value = '0.0₃0853'
$ws.Range("E21").Value = '  -6.77%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '237.03'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -1.03%  '

$ws.Range("E23").Value = '  -5.40%  '

$ws.Range("E24").Value = '  +0.12%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.34'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -5.65%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.27'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -4.13%  '

$ws.Range("E27").Value = '  -5.51%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '162.47'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -5.52%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '20.01'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.48%  '

$ws.Range("E30").Value = '  -3.46%  '

$ws.Range("E31").Value = '  -9.21%  '

$ws.Range("E32").Value = '  -0.66%  '

$ws.Range("E33").Value = '  -7.60%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0591'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -6.53%  '

$ws.Range("E35").Value = '  -0.01%  '

$ws.Range("E36").Value = '  +1.27%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0829'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -6.23%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.17'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -7.98%  '

$ws.Range("E39").Value = '  -8.73%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.81'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -7.57%  '

$ws.Range("E41").Value = '  -6.25%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.82'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -10.13%  '

$ws.Range("E43").Value = '  -5.58%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '93.10'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -7.91%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0899'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -11.79%  '

$ws.Range("D46").Value = '1.372.99'
$ws.Range("E46").Value = '  +4.65%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '15.61'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -10.96%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.28'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +4.07%  '

$ws.Range("E49").Value = '  -1.59%  '

$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").Value = '2.241.86'
$ws.Range("E50").Value = '  -0.43%  '

$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.25'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -8.04%  '
